$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 565, shifting rows 565:577 down to 569:581.
$ws.Range("A565:A568").EntireRow.Insert()

# New weekly data block (week of 2021-09-09, serial 44448), same product group
# (Hortaliza / Vega Modelo de Temuco / La Araucanía / Cebolla) as the rows around it.
$newRows = @(
    @{ Row = 565; H = "Morada(o)";        I = "1a (guarda)"; J = 250;  K = 10000; L = 10000; M = 10000; O = "Región de Arica y Parinacota"; P = 556 },
    @{ Row = 566; H = "Sin especificar";  I = "1a (guarda)"; J = 930;  K = 6000;  L = 6500;  M = 6204;  O = "Región de O'Higgins";           P = 345 },
    @{ Row = 567; H = "Sin especificar";  I = "1a (guarda)"; J = 890;  K = 5000;  L = 5000;  M = 5000;  O = "Región del Maule";              P = 278 },
    @{ Row = 568; H = "Sin especificar";  I = "Primera";     J = 650;  K = 4000;  L = 4000;  M = 4000;  O = "Perú";                          P = 222 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 10
    $ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($row, 3).Value = "La Araucanía"
    $ws.Cells.Item($row, 4).Value = 44448
    $ws.Cells.Item($row, 5).Value = 9
    $ws.Cells.Item($row, 6).Value = 100112004
    $ws.Cells.Item($row, 7).Value = "Cebolla"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "$/malla 18 kilos"
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
